$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The BOM row for the 3.3V LDO (row 7) switched parts: the old
# TLV70033DDCR / 296-27937-6-ND / "LDO 3.3V 200mA" is replaced with the
# newly released LP5907MFX-3.3/NOPB / 296-38557-6-ND /
# "Ultra-Low-Noise LDO 3.3V 250mA". Set the values left-to-right (E, H,
# then B) so new shared-string entries are appended in that order.
$ws.Range("E7").Value = "LP5907MFX-3.3/NOPB"
$ws.Range("H7").Value = "296-38557-6-ND"
$ws.Range("B7").Value = "Ultra-Low-Noise LDO 3.3V 250mA"

# Those same three cells also now use the formatting already applied to
# the "released"/highlighted row 10 instead of the plain row style. Copy
# that formatting (and only the formatting) from row 10 onto row 7.
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("H10").Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# The active selection moved from G18 to D16.
$ws.Range("D16").Select()
